# Swap the status/Interno Fiocruz/Externo à Fiocruz values between each
# pair of rows for the same docente, for the specific row-pairs affected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(5, 6),
    @(10, 11),
    @(26, 27),
    @(31, 32),
    @(35, 36),
    @(39, 40),
    @(41, 42)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Capture row1 values (columns B, C, D)
    $b1 = $ws.Cells.Item($r1, 2).Value2
    $c1 = $ws.Cells.Item($r1, 3).Value2
    $d1 = $ws.Cells.Item($r1, 4).Value2

    # Capture row2 values (columns B, C, D)
    $b2 = $ws.Cells.Item($r2, 2).Value2
    $c2 = $ws.Cells.Item($r2, 3).Value2
    $d2 = $ws.Cells.Item($r2, 4).Value2

    # Write row1 <- row2 values
    $ws.Cells.Item($r1, 2).Value = $b2
    $ws.Cells.Item($r1, 3).Value = $c2
    $ws.Cells.Item($r1, 4).Value = $d2

    # Write row2 <- row1 values
    $ws.Cells.Item($r2, 2).Value = $b1
    $ws.Cells.Item($r2, 3).Value = $c1
    $ws.Cells.Item($r2, 4).Value = $d1
}
